# edit.ps1 -- "fixed issues with workflow_2 graph"
#
# Applies the following content fixes to the single slide of the
# workflow_2 provenance graph deck:
#
#   1. Connector "Straight Arrow Connector 3" (shape Id=4, a bent
#      connector) gets its adj1 guide nudged from 50% to 54.609%.
#   2. Four command-line text boxes (shape Ids 96, 2, 89, 178) drop the
#      stray "preprocess.py " token from their run text and shrink to
#      match the now-shorter text.
#   3. The now-orphaned "Rectangle 48" shape (Id=49, the
#      "/home/pr/exp2/data.csv" label floating off the visible canvas)
#      is deleted.
#
# Notes on precision: Shape.Width/.Left/etc. are exposed as points
# (EMU / 12700) and round-trip through a lower-precision store
# internally, so naively assigning `targetEmu / 12700.0` can land the
# re-derived EMU one unit away from the intended target. To compensate,
# Set-WidthEmu/Set-LeftEmu nudge the assigned point value by tiny
# increments -- reading the property back (at full double precision via
# ToString("G17")) after every attempt -- until the value that would be
# re-derived on save matches the requested EMU exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

function Set-LeftEmu($shape, [long]$targetEmu) {
    $basePts = $targetEmu / 12700.0
    for ($k = 0; $k -lt 5000; $k++) {
        $candidate = $basePts + ($k * 0.0000001)
        $shape.Left = $candidate
        $readBack = [double]($shape.Left.ToString("G17"))
        $impliedEmu = [math]::Round($readBack * 12700.0)
        if ($impliedEmu -eq $targetEmu) {
            return
        }
    }
}

function Set-WidthEmu($shape, [long]$targetEmu) {
    $basePts = $targetEmu / 12700.0
    for ($k = 0; $k -lt 5000; $k++) {
        $candidate = $basePts + ($k * 0.0000001)
        $shape.Width = $candidate
        $readBack = [double]($shape.Width.ToString("G17"))
        $impliedEmu = [math]::Round($readBack * 12700.0)
        if ($impliedEmu -eq $targetEmu) {
            return
        }
    }
}

function Replace-RunText($textRange, [string]$oldText, [string]$newText) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -ge 0) {
        $sub = $textRange.Characters($idx + 1, $oldText.Length)
        $sub.Text = $newText
    }
}

# 1. Bent connector adj1: 50000 -> 54609 (out of 100000)
$connector = Get-ShapeById $s 4
if ($connector -ne $null) {
    $connector.Adjustments.Item(1) = 0.54609
}

# 2a. Shape Id=96 "Rectangle 95" -- drop "preprocess.py " and shrink.
$sh96 = Get-ShapeById $s 96
if ($sh96 -ne $null) {
    Replace-RunText $sh96.TextFrame.TextRange `
        "/bin/python3 preprocess.py train_model.py --" `
        "/bin/python3 train_model.py --"
    Set-WidthEmu $sh96 5577026
}

# 2b. Shape Id=2 "Rectangle 1" -- drop "preprocess.py " and shrink.
$sh2 = Get-ShapeById $s 2
if ($sh2 -ne $null) {
    Replace-RunText $sh2.TextFrame.TextRange `
        "preprocess.py train_model.py --" `
        "train_model.py --"
    Set-WidthEmu $sh2 5245957
}

# 2c. Shape Id=89 "Rectangle 88" -- drop "preprocess.py " and shrink.
$sh89 = Get-ShapeById $s 89
if ($sh89 -ne $null) {
    Replace-RunText $sh89.TextFrame.TextRange `
        "preprocess.py train_model.py --" `
        "train_model.py --"
    Set-WidthEmu $sh89 5370437
}

# 2d. Shape Id=178 "Rectangle 177" -- drop "preprocess.py ", nudge the
#     left edge by 1 EMU, and shrink.
$sh178 = Get-ShapeById $s 178
if ($sh178 -ne $null) {
    Replace-RunText $sh178.TextFrame.TextRange `
        "/bin/python3 preprocess.py evaluate_models.py " `
        "/bin/python3 evaluate_models.py "
    Set-LeftEmu $sh178 1062280
    Set-WidthEmu $sh178 4269138
}

# 3. Delete the stray "Rectangle 48" data.csv label (Id=49).
$sh49 = Get-ShapeById $s 49
if ($sh49 -ne $null) {
    $sh49.Delete()
}
